$wb = $excel.ActiveWorkbook

# 1. Delete the obsolete "Sheet" metadata row (row 16) from optimization_parameters.
#    This shifts the simulation_timepoints row up from 17 to 16, and drops the
#    now-unused "Sheet" shared string.
$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Activate()
$wsOpt.Rows.Item(16).Delete()
$wsOpt.Rows.Item(16).Select()

# 2. Refresh the recalculated weight on network_weights!D5 and move the
#    sheet's selection.
$wsWeights = $wb.Worksheets.Item("network_weights")
$wsWeights.Activate()
$wsWeights.Range("D5").Value = -0.97501548238480895
$wsWeights.Range("D9").Select()

# 3. Switch the active/selected sheet from optimization_parameters to threshold_b
$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
$wsThreshold.Range("A2").Select()
